# Updates cryptos list values (price/volume columns) and the two row swaps
# (NEARProtocol/Stellar and Stacks/Aave) per the target commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are assigned with a leading apostrophe so Excel stores them as literal
# text (matching the original inline-string cells) instead of re-parsing
# numeric-looking or percentage-looking strings as numbers.

$ws.Range("D2").Value = '''42.482.56'
$ws.Range("E2").Value = '''  -1.04%  '

$ws.Range("D3").Value = '''2.227.92'
$ws.Range("E3").Value = '''  -0.15%  '

$ws.Range("E4").Value = '''  -0.11%  '

$ws.Range("D5").Value = '''112.40'
$ws.Range("E5").Value = '''  -0.06%  '

$ws.Range("D6").Value = '''295.40'
$ws.Range("E6").Value = '''  +10.49%  '

$ws.Range("D7").Value = '''0.627'
$ws.Range("E7").Value = '''  +0.01%  '

$ws.Range("E8").Value = '''  -0.46%  '

$ws.Range("E9").Value = '''  +0.14%  '

$ws.Range("D10").Value = '''43.74'
$ws.Range("E10").Value = '''  -4.36%  '

$ws.Range("D11").Value = '''0.0917'
$ws.Range("E11").Value = '''  -0.89%  '

$ws.Range("D12").Value = '''54.51'
$ws.Range("E12").Value = '''  +1.12%  '

$ws.Range("D13").Value = '''8.64'
$ws.Range("E13").Value = '''  -5.12%  '

$ws.Range("D14").Value = '''1.05'
$ws.Range("E14").Value = '''  +20.84%  '

$ws.Range("E15").Value = '''  -1.04%  '

$ws.Range("D16").Value = '''14.98'
$ws.Range("E16").Value = '''  -1.63%  '

$ws.Range("D17").Value = '''2.564.49'
$ws.Range("E17").Value = '''  -0.16%  '

$ws.Range("D18").Value = '''2.229.99'
$ws.Range("E18").Value = '''  -0.63%  '

$ws.Range("D19").Value = '''42.505.69'
$ws.Range("E19").Value = '''  -1.59%  '

$ws.Range("D20").Value = '''7.20'
$ws.Range("E20").Value = '''  +7.76%  '

$ws.Range("E21").Value = '''  -1.30%  '

$ws.Range("D22").Value = '''73.58'
$ws.Range("E22").Value = '''  +2.84%  '

$ws.Range("E23").Value = '''  +15.80%  '

$ws.Range("D24").Value = '''2.39'
$ws.Range("E24").Value = '''  +1.79%  '

$ws.Range("D25").Value = '''236.30'
$ws.Range("E25").Value = '''  +2.54%  '

$ws.Range("D26").Value = '''8.88'
$ws.Range("E26").Value = '''  -3.48%  '

$ws.Range("E27").Value = '''  -1.66%  '

$ws.Range("D28").Value = '''11.45'
$ws.Range("E28").Value = '''  -5.74%  '

$ws.Range("D29").Value = '''2.21'
$ws.Range("E29").Value = '''  -0.89%  '

$ws.Range("D30").Value = '''175.51'
$ws.Range("E30").Value = '''  +1.44%  '

$ws.Range("D31").Value = '''37.37'
$ws.Range("E31").Value = '''  -7.15%  '

$ws.Range("E32").Value = '''  -4.48%  '

$ws.Range("E33").Value = '''  +1.70%  '

$ws.Range("D34").Value = '''0.0882'
$ws.Range("E34").Value = '''  -1.38%  '

$ws.Range("E35").Value = '''  +2.03%  '

$ws.Range("E36").Value = '''  +8.48%  '

$ws.Range("B37").Value = '''Stellar'
$ws.Range("C37").Value = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D37").Value = '''0.126'
$ws.Range("E37").Value = '''  -0.10%  '

$ws.Range("B38").Value = '''NEARProtocol'
$ws.Range("C38").Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").Value = '''4.19'
$ws.Range("E38").Value = '''  -0.07%  '

$ws.Range("E39").Value = '''  +0.99%  '

$ws.Range("E40").Value = '''  -2.01%  '

$ws.Range("E41").Value = '''  -4.42%  '

$ws.Range("D42").Value = '''72.08'
$ws.Range("E42").Value = '''  +1.86%  '

$ws.Range("D43").Value = '''0.230'
$ws.Range("E43").Value = '''  -0.16%  '

$ws.Range("E44").Value = '''  -0.19%  '

$ws.Range("D45").Value = '''12.31'
$ws.Range("E45").Value = '''  -7.00%  '

$ws.Range("E46").Value = '''  -1.23%  '

$ws.Range("D47").Value = '''5.38'
$ws.Range("E47").Value = '''  -4.65%  '

$ws.Range("E48").Value = '''  +4.05%  '

$ws.Range("B49").Value = '''Aave'
$ws.Range("C49").Value = '''https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '''102.04'
$ws.Range("E49").Value = '''  +2.52%  '

$ws.Range("B50").Value = '''Stacks'
$ws.Range("C50").Value = '''https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").Value = '''1.65'
$ws.Range("E50").Value = '''  +7.32%  '

$ws.Range("E51").Value = '''  +0.92%  '
